$d = $word.ActiveDocument

$replacements = @(
    @("460÷2=230, 0", "623÷9=69, 2"),
    @("642÷6=107, 0", "522÷7=74, 4"),
    @("788÷3=262, 2", "549÷7=78, 3"),
    @("114÷9=12, 6", "465÷7=66, 3"),
    @("432÷3=144, 0", "451÷2=225, 1"),
    @("307÷8=38, 3", "237÷2=118, 1"),
    @("493÷3=164, 1", "731÷8=91, 3"),
    @("469÷3=156, 1", "485÷5=97, 0"),
    @("696÷5=139, 1", "259÷5=51, 4"),
    @("176÷4=44, 0", "277÷9=30, 7"),
    @("388÷9=43, 1", "327÷9=36, 3"),
    @("467÷7=66, 5", "765÷6=127, 3"),
    @("799÷8=99, 7", "128÷2=64, 0"),
    @("836÷3=278, 2", "347÷9=38, 5"),
    @("351÷2=175, 1", "571÷7=81, 4"),
    @("224÷5=44, 4", "627÷4=156, 3"),
    @("298÷6=49, 4", "337÷9=37, 4"),
    @("769÷4=192, 1", "928÷4=232, 0"),
    @("293÷9=32, 5", "231÷6=38, 3"),
    @("116÷2=58, 0", "587÷8=73, 3"),
    @("642÷8=80, 2", "549÷9=61, 0"),
    @("194÷3=64, 2", "756÷5=151, 1"),
    @("227÷6=37, 5", "612÷5=122, 2"),
    @("973÷4=243, 1", "910÷6=151, 4"),
    @("454÷5=90, 4", "447÷8=55, 7")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

Write-Host "Done applying $($replacements.Count) replacements"
